# Applies the crypto price/volume refresh described in the commit:
# "Updated cryptos list on Sun Jul 30 07:30:57 UTC 2023 with GitHub Actions"
#
# The sheet stores Price (D) and Volume(1h) (E) columns as plain text
# (e.g. "244.00", "1.262.09", "  +0.09%  "). Setting .Value directly would
# let Excel's auto-detection coerce these into numbers (dropping trailing
# zeros / using scientific notation for tiny values), so we force the
# cell's NumberFormat to Text ("@") before writing the new literal string.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Cell = 'D2'; Value = '29.337.01' }
    @{ Cell = 'E2'; Value = '  -0.02%  ' }
    @{ Cell = 'D3'; Value = '1.877.92' }
    @{ Cell = 'E3'; Value = '  +0.24%  ' }
    @{ Cell = 'D5'; Value = '0.7107' }
    @{ Cell = 'E5'; Value = '  -0.27%  ' }
    @{ Cell = 'D6'; Value = '242.35' }
    @{ Cell = 'E6'; Value = '  +0.25%  ' }
    @{ Cell = 'E7'; Value = '  +0.18%  ' }
    @{ Cell = 'D8'; Value = '0.08002' }
    @{ Cell = 'E8'; Value = '  +2.65%  ' }
    @{ Cell = 'E9'; Value = '  +1.54%  ' }
    @{ Cell = 'D10'; Value = '24.96' }
    @{ Cell = 'E10'; Value = '  -0.57%  ' }
    @{ Cell = 'D11'; Value = '0.08297' }
    @{ Cell = 'E11'; Value = '  -1.46%  ' }
    @{ Cell = 'D12'; Value = '1.889.97' }
    @{ Cell = 'E12'; Value = '  +0.94%  ' }
    @{ Cell = 'D13'; Value = '5.249' }
    @{ Cell = 'E13'; Value = '  +0.31%  ' }
    @{ Cell = 'D14'; Value = '94.41' }
    @{ Cell = 'E14'; Value = '  +3.66%  ' }
    @{ Cell = 'D15'; Value = '0.7146' }
    @{ Cell = 'E15'; Value = '  +0.34%  ' }
    @{ Cell = 'D16'; Value = '6.362' }
    @{ Cell = 'E16'; Value = '  +4.40%  ' }
    @{ Cell = 'D17'; Value = '0.000008523' }
    @{ Cell = 'E17'; Value = '  +3.25%  ' }
    @{ Cell = 'D18'; Value = '29.378.71' }
    @{ Cell = 'E18'; Value = '  +0.09%  ' }
    @{ Cell = 'D19'; Value = '244.00' }
    @{ Cell = 'E19'; Value = '  +1.66%  ' }
    @{ Cell = 'B20'; Value = 'WrappedliquidstakedEther2.0' }
    @{ Cell = 'C20'; Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth' }
    @{ Cell = 'D20'; Value = '2.165.18' }
    @{ Cell = 'E20'; Value = '  +2.82%  ' }
    @{ Cell = 'B21'; Value = 'Avalanche' }
    @{ Cell = 'C21'; Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax' }
    @{ Cell = 'D21'; Value = '13.28' }
    @{ Cell = 'E21'; Value = '  +0.59%  ' }
    @{ Cell = 'E22'; Value = '  +0.24%  ' }
    @{ Cell = 'D23'; Value = '7.802' }
    @{ Cell = 'E23'; Value = '  +0.68%  ' }
    @{ Cell = 'D24'; Value = '1.002' }
    @{ Cell = 'E24'; Value = '  +0.15%  ' }
    @{ Cell = 'E25'; Value = '  -1.79%  ' }
    @{ Cell = 'D26'; Value = '9.071' }
    @{ Cell = 'E26'; Value = '  +0.31%  ' }
    @{ Cell = 'D27'; Value = '162.74' }
    @{ Cell = 'E27'; Value = '  -0.21%  ' }
    @{ Cell = 'D28'; Value = '18.54' }
    @{ Cell = 'E28'; Value = '  +0.17%  ' }
    @{ Cell = 'D29'; Value = '1.509' }
    @{ Cell = 'E29'; Value = '  -0.11%  ' }
    @{ Cell = 'D30'; Value = '4.417' }
    @{ Cell = 'E30'; Value = '  -0.03%  ' }
    @{ Cell = 'D31'; Value = '4.324' }
    @{ Cell = 'E31'; Value = '  +0.08%  ' }
    @{ Cell = 'E32'; Value = '  -7.85%  ' }
    @{ Cell = 'D33'; Value = '0.05376' }
    @{ Cell = 'E33'; Value = '  +1.69%  ' }
    @{ Cell = 'D34'; Value = '1.936' }
    @{ Cell = 'E34'; Value = '  -0.02%  ' }
    @{ Cell = 'D35'; Value = '0.7673' }
    @{ Cell = 'E35'; Value = '  +3.67%  ' }
    @{ Cell = 'D36'; Value = '1.185' }
    @{ Cell = 'E36'; Value = '  +0.60%  ' }
    @{ Cell = 'D37'; Value = '2.688' }
    @{ Cell = 'E37'; Value = '  -0.44%  ' }
    @{ Cell = 'D38'; Value = '0.01886' }
    @{ Cell = 'E38'; Value = '  +0.54%  ' }
    @{ Cell = 'D39'; Value = '1.262.09' }
    @{ Cell = 'E39'; Value = '  +3.92%  ' }
    @{ Cell = 'D40'; Value = '2.753' }
    @{ Cell = 'E40'; Value = '  +0.88%  ' }
    @{ Cell = 'D41'; Value = '6.515' }
    @{ Cell = 'E41'; Value = '  -0.43%  ' }
    @{ Cell = 'D42'; Value = '113.12' }
    @{ Cell = 'B43'; Value = 'TrustWalletToken' }
    @{ Cell = 'C43'; Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt' }
    @{ Cell = 'D43'; Value = '0.9048' }
    @{ Cell = 'E43'; Value = '  +1.98%  ' }
    @{ Cell = 'B44'; Value = 'Aave' }
    @{ Cell = 'C44'; Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave' }
    @{ Cell = 'D44'; Value = '74.27' }
    @{ Cell = 'E44'; Value = '  +1.89%  ' }
    @{ Cell = 'D45'; Value = '0.00000000133' }
    @{ Cell = 'E45'; Value = '  +8.55%  ' }
    @{ Cell = 'E46'; Value = '  +0.18%  ' }
    @{ Cell = 'D47'; Value = '2.027.54' }
    @{ Cell = 'E47'; Value = '  +0.48%  ' }
    @{ Cell = 'E48'; Value = '  +0.41%  ' }
    @{ Cell = 'D49'; Value = '1.799' }
    @{ Cell = 'E49'; Value = '  -0.23%  ' }
    @{ Cell = 'E50'; Value = '  +0.57%  ' }
    @{ Cell = 'D51'; Value = '0.4373' }
    @{ Cell = 'E51'; Value = '  +1.29%  ' }
)

foreach ($u in $updates) {
    $col = $u.Cell -replace '[0-9]+$', ''
    $rng = $ws.Range($u.Cell)
    if ($col -eq 'D' -or $col -eq 'E') {
        # Keep these text-typed so values like "244.00" or "0.08002"
        # round-trip exactly instead of becoming numeric.
        $rng.NumberFormat = "@"
    }
    $rng.Value = $u.Value
}

